# "Fruta / hortaliza, semanal" weekly update:
# A new weekly price-report row for Mango (Terminal Hortofrutícola Agro
# Chillán) is inserted at row 62, pushing the previously-existing rows
# 62-122 down to 63-123 (dimension grows from A1:T122 to A1:T123).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above the current row 62; Excel shifts rows
# 62-122 down to 63-123 and carries their formatting with them.
$ws.Rows.Item(62).Insert()

# Populate the newly inserted row 62 with this week's record.
$ws.Cells.Item(62, 1).Value2  = 7
$ws.Cells.Item(62, 2).Value2  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(62, 3).Value2  = "Ñuble"
$ws.Cells.Item(62, 4).Value2  = 45040
$ws.Cells.Item(62, 5).Value2  = 16
$ws.Cells.Item(62, 6).Value2  = "Fruta"
$ws.Cells.Item(62, 7).Value2  = 100108
$ws.Cells.Item(62, 8).Value2  = "Tropicales y subtropicales"
$ws.Cells.Item(62, 9).Value2  = 100108002
$ws.Cells.Item(62, 10).Value2 = "Mango"
$ws.Cells.Item(62, 11).Value2 = "Sin especificar"
$ws.Cells.Item(62, 12).Value2 = "Primera"
$ws.Cells.Item(62, 13).Value2 = 90
$ws.Cells.Item(62, 14).Value2 = 7000
$ws.Cells.Item(62, 15).Value2 = 8000
$ws.Cells.Item(62, 16).Value2 = 7556
$ws.Cells.Item(62, 17).Value2 = "$/bandeja 4 kilos"
$ws.Cells.Item(62, 18).Value2 = "Perú"
$ws.Cells.Item(62, 19).Value2 = 1889
$ws.Cells.Item(62, 20).Value2 = 4
